# Weekly update: insert a new Mango price record for Terminal Hortofrutícola
# Agro Chillán (Fruta / hortaliza, semanal).
#
# The new observation is inserted as row 164, pushing the previously
# existing rows 164-187 down to 165-188 (and updating the sheet dimension
# from A1:T187 to A1:T188 automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164 - shifts rows 164:187 down to 165:188
# and inherits the formatting (incl. the date number format in column D)
# from the row above.
$ws.Rows(164).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(164, 1).Value  = 7
$ws.Cells.Item(164, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(164, 3).Value  = "Ñuble"
$ws.Cells.Item(164, 4).Value  = 45173
$ws.Cells.Item(164, 5).Value  = 16
$ws.Cells.Item(164, 6).Value  = "Fruta"
$ws.Cells.Item(164, 7).Value  = 100108
$ws.Cells.Item(164, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value  = 100108002
$ws.Cells.Item(164, 10).Value = "Mango"
$ws.Cells.Item(164, 11).Value = "Sin especificar"
$ws.Cells.Item(164, 12).Value = "Primera"
$ws.Cells.Item(164, 13).Value = 60
$ws.Cells.Item(164, 14).Value = 11000
$ws.Cells.Item(164, 15).Value = 11000
$ws.Cells.Item(164, 16).Value = 11000
$ws.Cells.Item(164, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(164, 18).Value = "Brasil"
$ws.Cells.Item(164, 19).Value = 2750
$ws.Cells.Item(164, 20).Value = 4
